$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.043.26'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '2.360.30'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'0.688"
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("D6").Value = "'240.13"
$ws.Range("E6").Value = '  +2.07%  '
$ws.Range("D7").Value = "'76.03"
$ws.Range("E7").Value = '  +4.68%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +14.18%  '
$ws.Range("E10").Value = '  +2.38%  '
$ws.Range("D11").Value = "'57.19"
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").Value = "'33.18"
$ws.Range("E12").Value = '  +17.86%  '
$ws.Range("E13").Value = '  +11.36%  '
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '2.707.86'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = "'16.64"
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").Value = "'0.914"
$ws.Range("E17").Value = '  +3.66%  '
$ws.Range("D18").Value = '2.360.60'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '43.936.29'
$ws.Range("E20").Value = '  +1.59%  '
$ws.Range("E21").Value = '  +5.15%  '
$ws.Range("D22").Value = "'77.52"
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("D23").Value = "'259.08"
$ws.Range("E23").Value = '  +2.91%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = "'3.71"
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("D26").Value = "'1.84"
$ws.Range("E26").Value = '  +17.56%  '
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = "'10.81"
$ws.Range("E28").Value = '  +3.50%  '
$ws.Range("E29").Value = '  +2.12%  '
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("D31").Value = "'174.97"
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("E33").Value = '  +4.13%  '
$ws.Range("D34").Value = "'0.0758"
$ws.Range("E34").Value = '  +6.59%  '
$ws.Range("E35").Value = '  +2.49%  '
$ws.Range("E36").Value = '  +3.59%  '
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("D40").Value = "'0.0283"
$ws.Range("E40").Value = '  +3.66%  '
$ws.Range("D41").Value = "'0.214"
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = "'19.44"
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'9.23"
$ws.Range("E43").Value = '  +3.42%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = "'0.108"
$ws.Range("E44").Value = '  +11.66%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = "'4.77"
$ws.Range("E46").Value = '  +7.70%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'2.58"
$ws.Range("E47").Value = '  +11.46%  '
$ws.Range("E48").Value = '  +3.29%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = "'57.02"
$ws.Range("E50").Value = '  +8.95%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'100.57"
$ws.Range("E51").Value = '  +2.79%  '
